# Apply the "remove RG" edit to the certificate model (slide 1, shape "Rectangle 5").
#
# 1. Rewrite "portador do RG nº {{RG}} e CPF nº " -> "portador do CPF nº "
#    (drop the RG mention entirely, keep the CPF one).
# 2. Merge "São Carlos" + ", " into a single run "São Carlos, ".
# 3. Shrink the text placeholder shape - its text got shorter once the RG
#    sentence was trimmed, so PowerPoint's autofit reduced its height.
#    (Done last: the shape has <a:spAutoFit/>, so the host recomputes the
#    height automatically whenever the text changes - setting it first
#    would just get overwritten by the following edits.)

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item("Rectangle 5")   # the certificate body text box
$tr = $sh.TextFrame.TextRange

# --- 1. Drop the RG reference, keep the CPF one ----------------------------

# "portador do RG nº " -> delete "RG nº " -> "portador do "
$txt = $tr.Text
$idx = $txt.IndexOf("RG n")
$rng = $tr.Characters($idx + 1, "RG nº ".Length)
[void]$rng.Delete()

# Split "portador do " into separate "portador " / "do " runs
$txt = $tr.Text
$idx = $txt.IndexOf("portador do ")
$doStart = $idx + "portador ".Length
$rng = $tr.Characters($doStart + 1, "do ".Length)
$rng.Text = "do "

# "{{RG}} " -> "CPF "
$txt = $tr.Text
$idx = $txt.IndexOf("{{RG}} ")
$rng = $tr.Characters($idx + 1, "{{RG}} ".Length)
$rng.Text = "CPF "

# "e CPF nº " -> delete "e CPF " -> "nº "
$txt = $tr.Text
$idx = $txt.IndexOf("e CPF n")
$rng = $tr.Characters($idx + 1, "e CPF ".Length)
[void]$rng.Delete()

# --- 2. Merge "São Carlos" and ", " into a single run ----------------------

$txt = $tr.Text
$idx = $txt.IndexOf("Carlos")
$scStart = $idx - "São ".Length
$rng = $tr.Characters($scStart + 1, "São Carlos".Length)
[void]$rng.InsertAfter(", ")

$txt = $tr.Text
$idx = $txt.IndexOf("São Carlos, ")
$commaStart = $idx + "São Carlos, ".Length
$rng = $tr.Characters($commaStart + 1, 2)
[void]$rng.Delete()

# --- 3. Resize the shape (height only; width/position stay the same) ------

$sh.Height = 2462213 / 12700
